# Insert a new data row at row 217 (pushing the existing rows 217-301 down
# to 218-302), then populate the newly inserted row with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(217).Insert()

$ws.Range("A217").Value = 5
$ws.Range("B217").Value = "Macroferia Regional de Talca"
$ws.Range("C217").Value = "Maule"
$ws.Range("D217").Value = 44559
$ws.Range("E217").Value = 7
$ws.Range("F217").Value = "Fruta"
$ws.Range("G217").Value = 100109
$ws.Range("H217").Value = "Uva"
$ws.Range("I217").Value = 100109001
$ws.Range("J217").Value = "Uva"
$ws.Range("K217").Value = "Superior Seedless"
$ws.Range("L217").Value = "Segunda"
$ws.Range("M217").Value = 300
$ws.Range("N217").Value = 11000
$ws.Range("O217").Value = 11000
$ws.Range("P217").Value = 11000
$ws.Range("Q217").Value = "`$/bandeja 10 kilos"
$ws.Range("R217").Value = "Provincia de Limarí"
$ws.Range("S217").Value = 1100
$ws.Range("T217").Value = 10
